$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the selected/active cell shown in the sheet view
$ws.Range("J20").Select()

# Row 13: invests (L) for 2021
$ws.Range("L13").Value = 1384.127

# Row 14: beforeschool (H) and invests (L) for 2022
$ws.Range("H14").Value = 2544
$ws.Range("H14").HorizontalAlignment = -4108
$ws.Range("L14").Value = 863.57799999999997

# Row 15: beforeschool (H) for 2023 (cell already existed, keep its style)
$ws.Range("H15").Value = 2139

# Row 16: invests (L) for 2021
$ws.Range("L16").Value = 15703.486000000001

# Row 17: invests (L) for 2021
$ws.Range("L17").Value = 10326.647999999999

# Row 18: beforeschool (H) and invests (L) for 2022 (cells already existed, keep style)
$ws.Range("H18").Value = 8508
$ws.Range("L18").Value = 26447.994999999999

# Row 19: beforeschool (H) and invests (L) for 2022 (cells already existed, keep style)
$ws.Range("H19").Value = 5843
$ws.Range("L19").Value = 7713.7349999999997

# Row 20: beforeschool (H) for 2023
$ws.Range("H20").Value = 8508
$ws.Range("H20").HorizontalAlignment = -4108

# Row 21: beforeschool (H) for 2023
$ws.Range("H21").Value = 5703
$ws.Range("H21").HorizontalAlignment = -4108
